$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 14:22"

# Swap the country names for rows 61 (was Grecia) and 62 (was Kuwait)
# so that row 61 now shows Kuwait and row 62 now shows Grecia.
$ws.Range("A61").Value = "Kuwait"
$ws.Range("A62").Value = "Grecia"

# Row 17 (Iran) - updated covid figures
$ws.Range("B17").Value = 36535
$ws.Range("C17").Value = 806
$ws.Range("E17").Value = 31996
$ws.Range("G17").Value = 112
$ws.Range("H17").Value = 4289

# Row 23 (Peru) - updated covid figures
$ws.Range("B23").Value = 17567
$ws.Range("C23").Value = 812
$ws.Range("E23").Value = 14865
$ws.Range("G23").Value = 131
$ws.Range("H23").Value = 2152

# Row 39 - updated covid figures
$ws.Range("D39").Value = 5526
$ws.Range("E39").Value = 2281
$ws.Range("F39").Value = 69
$ws.Range("G39").Value = 9
$ws.Range("H39").Value = 403

# Row 61 (now Kuwait) - updated covid figures
$ws.Range("B61").Value = 2614
$ws.Range("C61").Value = 215
$ws.Range("D61").Value = 613
$ws.Range("E61").Value = 1986
$ws.Range("F61").Value = 60
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 15

# Row 62 (now Grecia) - updated covid figures
$ws.Range("B62").Value = 2463
$ws.Range("D62").Value = 577
$ws.Range("E62").Value = 1759
$ws.Range("F62").Value = 52
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 127

# Row 73 - updated covid figures
$ws.Range("B73").Value = 1592
$ws.Range("C73").Value = 44
$ws.Range("D73").Value = 1013
$ws.Range("E73").Value = 558
$ws.Range("G73").Value = 1
$ws.Range("H73").Value = 21

# Row 114 - updated covid figures
$ws.Range("B114").Value = 379
$ws.Range("C114").Value = 11
$ws.Range("E114").Value = 265

# Row 124 - updated covid figures
$ws.Range("B124").Value = 270
$ws.Range("C124").Value = 2
$ws.Range("E124").Value = 45
